$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a plain number (e.g. "589.38") must be forced
# to stay as text, matching the original inlineStr cells. Prefixing with an
# apostrophe forces text entry; resetting the style back to Normal afterwards
# avoids leaving a stray quote-prefix number format on the cell.
$textCells = @(
    "D5","D6","D9","D11","D13","D14","D20","D21","D22","D24","D25",
    "D32","D33","D36","D41","D42","D48","D49","D50"
)

$updates = @(
    @{ Cell = "D2";  Value = "62.785.92" }
    @{ Cell = "E2";  Value = "  -1.27%  " }
    @{ Cell = "D3";  Value = "3.172.55" }
    @{ Cell = "E3";  Value = "  -4.19%  " }
    @{ Cell = "E4";  Value = "  +0.02%  " }
    @{ Cell = "D5";  Value = "589.38" }
    @{ Cell = "E5";  Value = "  -2.83%  " }
    @{ Cell = "D6";  Value = "136.75" }
    @{ Cell = "E6";  Value = "  -3.86%  " }
    @{ Cell = "E7";  Value = "  +0.03%  " }
    @{ Cell = "D8";  Value = "3.169.39" }
    @{ Cell = "E8";  Value = "  -4.35%  " }
    @{ Cell = "D9";  Value = "0.512" }
    @{ Cell = "E9";  Value = "  -1.51%  " }
    @{ Cell = "E10"; Value = "  -4.91%  " }
    @{ Cell = "D11"; Value = "5.28" }
    @{ Cell = "E11"; Value = "  -3.72%  " }
    @{ Cell = "E12"; Value = "  -2.61%  " }
    @{ Cell = "D13"; Value = "0.0000236" }
    @{ Cell = "E13"; Value = "  -4.90%  " }
    @{ Cell = "D14"; Value = "34.93" }
    @{ Cell = "E14"; Value = "  +0.76%  " }
    @{ Cell = "D15"; Value = "3.694.55" }
    @{ Cell = "E15"; Value = "  -4.17%  " }
    @{ Cell = "E16"; Value = "  -1.86%  " }
    @{ Cell = "D17"; Value = "3.171.51" }
    @{ Cell = "E17"; Value = "  -4.15%  " }
    @{ Cell = "D18"; Value = "62.771.26" }
    @{ Cell = "E18"; Value = "  -1.45%  " }
    @{ Cell = "E19"; Value = "  -3.85%  " }
    @{ Cell = "D20"; Value = "457.73" }
    @{ Cell = "E20"; Value = "  -4.75%  " }
    @{ Cell = "D21"; Value = "13.95" }
    @{ Cell = "E21"; Value = "  -1.48%  " }
    @{ Cell = "D22"; Value = "0.711" }
    @{ Cell = "E22"; Value = "  -3.04%  " }
    @{ Cell = "E23"; Value = "  -5.95%  " }
    @{ Cell = "D24"; Value = "13.41" }
    @{ Cell = "E24"; Value = "  -2.19%  " }
    @{ Cell = "D25"; Value = "83.35" }
    @{ Cell = "E25"; Value = "  -1.79%  " }
    @{ Cell = "E26"; Value = "  +0.00%  " }
    @{ Cell = "E27"; Value = "  -3.19%  " }
    @{ Cell = "E28"; Value = "  -0.09%  " }
    @{ Cell = "E29"; Value = "  -6.87%  " }
    @{ Cell = "E30"; Value = "  -4.59%  " }
    @{ Cell = "E31"; Value = "  -5.87%  " }
    @{ Cell = "D32"; Value = "27.36" }
    @{ Cell = "E32"; Value = "  -5.37%  " }
    @{ Cell = "D33"; Value = "0.103" }
    @{ Cell = "E33"; Value = "  -2.85%  " }
    @{ Cell = "E34"; Value = "  -6.17%  " }
    @{ Cell = "E35"; Value = "  -6.44%  " }
    @{ Cell = "D36"; Value = "5.84" }
    @{ Cell = "E36"; Value = "  -3.18%  " }
    @{ Cell = "E37"; Value = "  -3.27%  " }
    @{ Cell = "D38"; Value = "0.0₃0706" }
    @{ Cell = "E38"; Value = "  -5.16%  " }
    @{ Cell = "E39"; Value = "  -3.31%  " }
    @{ Cell = "E40"; Value = "  -2.10%  " }
    @{ Cell = "D41"; Value = "401.50" }
    @{ Cell = "E41"; Value = "  -7.04%  " }
    @{ Cell = "D42"; Value = "8.11" }
    @{ Cell = "E42"; Value = "  -2.80%  " }
    @{ Cell = "E43"; Value = "  -3.28%  " }
    @{ Cell = "D44"; Value = "2.777.89" }
    @{ Cell = "E44"; Value = "  -9.79%  " }
    @{ Cell = "E45"; Value = "  -4.28%  " }
    @{ Cell = "E46"; Value = "  -3.16%  " }
    @{ Cell = "D48"; Value = "125.22" }
    @{ Cell = "E48"; Value = "  +0.75%  " }
    @{ Cell = "D49"; Value = "25.59" }
    @{ Cell = "E49"; Value = "  -2.97%  " }
    @{ Cell = "D50"; Value = "34.78" }
    @{ Cell = "E50"; Value = "  -3.81%  " }
)

foreach ($u in $updates) {
    $cell = $u.Cell
    $value = $u.Value
    if ($textCells -contains $cell) {
        # Force text entry so a numeric-looking value (e.g. "589.38") does
        # not get auto-converted into a real number.
        $ws.Range($cell).Value = "'" + $value
        $ws.Range($cell).Style = "Normal"
    } else {
        $ws.Range($cell).Value = $value
    }
}
